$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value2 = '30.515.71'
$ws.Range('E2').Value2 = '  +2.77%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value2 = '2.122.11'
$ws.Range('E3').Value2 = '  +1.65%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value2 = '1.008'
$ws.Range('E4').Value2 = '  -0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '347.14'
$ws.Range('E5').Value2 = '  +0.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '1.007'
$ws.Range('E6').Value2 = '  -0.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value2 = '0.5234'
$ws.Range('E7').Value2 = '  +1.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '0.4481'
$ws.Range('E8').Value2 = '  +1.80%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '54.39'
$ws.Range('E9').Value2 = '  +5.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '0.09404'
$ws.Range('E10').Value2 = '  +1.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value2 = '1.181'
$ws.Range('E11').Value2 = '  +0.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value2 = '25.27'
$ws.Range('E12').Value2 = '  -0.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value2 = '8.723'
$ws.Range('E13').Value2 = '  +7.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '6.975'
$ws.Range('E14').Value2 = '  +3.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '2.100.98'
$ws.Range('E15').Value2 = '  +0.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value2 = '102.38'
$ws.Range('E16').Value2 = '  +2.86%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value2 = '0.00001168'
$ws.Range('E17').Value2 = '  +0.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value2 = '1.008'
$ws.Range('E18').Value2 = '  -0.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '21.47'
$ws.Range('E19').Value2 = '  +2.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '0.06727'
$ws.Range('E20').Value2 = '  +0.81%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '6.352'
$ws.Range('E21').Value2 = '  +2.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '1.007'
$ws.Range('E22').Value2 = '  -0.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value2 = '30.553.20'
$ws.Range('E23').Value2 = '  +2.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value2 = '12.77'
$ws.Range('E24').Value2 = '  +0.44%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '2.331'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value2 = '2.363.95'
$ws.Range('E26').Value2 = '  +1.03%  '
$ws.Range('E27').Value2 = '  +1.70%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value2 = '2.553'
$ws.Range('E28').Value2 = '  +0.98%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value2 = '162.63'
$ws.Range('E29').Value2 = '  -0.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '134.45'
$ws.Range('E30').Value2 = '  +1.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '1.163'
$ws.Range('E31').Value2 = '  +1.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '1.778'
$ws.Range('E32').Value2 = '  +9.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '0.1063'
$ws.Range('E33').Value2 = '  +0.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '6.906'
$ws.Range('E34').Value2 = '  +11.86%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '6.304'
$ws.Range('E35').Value2 = '  +1.64%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '3.964'
$ws.Range('E36').Value2 = '  +0.35%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '10.63'
$ws.Range('E37').Value2 = '  +4.81%  '
$ws.Range('E38').Value2 = '  +3.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value2 = '0.06876'
$ws.Range('E39').Value2 = '  +2.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '0.7141'
$ws.Range('E40').Value2 = '  +3.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '12.67'
$ws.Range('E41').Value2 = '  +2.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '0.2251'
$ws.Range('E42').Value2 = '  -1.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '1.337'
$ws.Range('E43').Value2 = '  +4.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '0.6943'
$ws.Range('E44').Value2 = '  +4.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '14.70'
$ws.Range('E45').Value2 = '  +4.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '2.398'
$ws.Range('E46').Value2 = '  +4.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value2 = '1.006'
$ws.Range('E47').Value2 = '  -0.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '1.331'
$ws.Range('E48').Value2 = '  +14.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '3.660'
$ws.Range('E49').Value2 = '  +0.94%  '
$ws.Range('B50').Value2 = 'BabyDogeCoin'
$ws.Range('C50').Value2 = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '0.00000000349'
$ws.Range('E50').Value2 = '  +3.14%  '
$ws.Range('B51').Value2 = 'EOS'
$ws.Range('C51').Value2 = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '1.230'
$ws.Range('E51').Value2 = '  +1.03%  '
